$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 55.3718303775116
$ws.Range("E2").Value = 56.7933697165468
$ws.Range("K2").Value = 54.8090076208069
$ws.Range("L2").Value = 49.0282605200155
$ws.Range("N2").Value = 47.4912333237318

$ws.Range("B3").Value = 40.5493307668479
$ws.Range("K3").Value = 36.9055302021312
$ws.Range("L3").Value = 36.9057907051389
$ws.Range("N3").Value = 41.6894250824717

$ws.Range("B4").Value = 34.7872415482579
$ws.Range("K4").Value = 30.4357942185907
$ws.Range("N4").Value = 43.0029533260978

$ws.Range("B5").Value = 64.7240480131028
$ws.Range("K5").Value = 66.5933688883394
$ws.Range("L5").Value = 57.2690723237937
$ws.Range("N5").Value = 60.3471262597791

$ws.Range("B6").Value = 67.2054760884641
$ws.Range("C6").Value = 73.150166637206
$ws.Range("K6").Value = 67.292535125419
$ws.Range("L6").Value = 63.9033952437169
$ws.Range("N6").Value = 58.9468339412676

$ws.Range("B7").Value = 68.518253119246
$ws.Range("K7").Value = 70.8149579691965
$ws.Range("L7").Value = 60.4866837966034
$ws.Range("N7").Value = 63.6828257869627

$ws.Range("B8").Value = 61.895417279981
$ws.Range("K8").Value = 63.3861854984241
$ws.Range("L8").Value = 54.11746546684
$ws.Range("N8").Value = 56.0767240057917
